# Generate Report for Handoff
# Adds a new handed-off file (d719813f-901f-467c-8dbc-a1e554fac620.md) as
# row 3 on the "Overview", "zh-cn" and "de-de" worksheets, mirroring the
# existing 6be85808-... row, and expands each table/dimension accordingly.

$wb = $excel.ActiveWorkbook

$newBase = "d719813f-901f-467c-8dbc-a1e554fac620"
$newMd = "$newBase.md"
$newMdDisplayOverview = "e2e\$newBase.md"

$commitSha = "c35798f56918125b7f7928ac405e9a5c99e1139e"
$zhXlf = "$newBase.$commitSha.zh-cn.xlf"
$deXlf = "$newBase.$commitSha.de-de.xlf"

$zhDateTime = "2016-08-20 14:43:43"
$deDateTime = "2016-08-20 14:43:47"
$epoch = "0001-01-01 00:00:00"

$baseHyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8e0a7de394663b947ee039cea9d61fac63d1727/e2e/$newMd"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Cells.Item(3, 1).Value = $newMd
$wsOverview.Cells.Item(3, 3).Value = ".md"
$wsOverview.Cells.Item(3, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item(3, 7).Value = $deDateTime
$wsOverview.Cells.Item(3, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $baseHyperlinkUrl, $null, $null, $newMdDisplayOverview) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Cells.Item(3, 2).Value = ".md"
$wsZhCn.Cells.Item(3, 3).Value = "Ready for handoff"
$wsZhCn.Cells.Item(3, 4).Value = "e2e"
$wsZhCn.Cells.Item(3, 5).Value = "ht"
$wsZhCn.Cells.Item(3, 6).Value = "'False"
$wsZhCn.Cells.Item(3, 7).Value = $zhXlf
$wsZhCn.Cells.Item(3, 8).Value = $zhDateTime
$wsZhCn.Cells.Item(3, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(3, 11).Value = $epoch
$wsZhCn.Cells.Item(3, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Cells.Item(3, 13).Value = "'True"
$wsZhCn.Cells.Item(3, 15).Value = "'False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $baseHyperlinkUrl, $null, $null, $newMd) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Cells.Item(3, 2).Value = ".md"
$wsDeDe.Cells.Item(3, 3).Value = "Ready for handoff"
$wsDeDe.Cells.Item(3, 4).Value = "e2e"
$wsDeDe.Cells.Item(3, 5).Value = "ht"
$wsDeDe.Cells.Item(3, 6).Value = "'False"
$wsDeDe.Cells.Item(3, 7).Value = $deXlf
$wsDeDe.Cells.Item(3, 8).Value = $deDateTime
$wsDeDe.Cells.Item(3, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(3, 11).Value = $epoch
$wsDeDe.Cells.Item(3, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Cells.Item(3, 13).Value = "'True"
$wsDeDe.Cells.Item(3, 15).Value = "'False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $baseHyperlinkUrl, $null, $null, $newMd) | Out-Null
